# Updates the "Price" (D) and "Volume(1h)" (E) columns on the active sheet
# to reflect the refreshed cryptocurrency symbol list, as produced by the
# "Updated symbol list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. All of these cells hold text (the
# workbook stores prices as strings, not numbers), so we force a text
# number format before writing the value to keep Excel from re-interpreting
# the string as a numeric value.
$updates = [ordered]@{
    "D2"  = "242.77"
    "D3"  = "23.13"
    "D4"  = "5.387"
    "D5"  = "0.05981"
    "D6"  = "3.401"
    "D7"  = "6.491"
    "D8"  = "0.8126"
    "D9"  = "0.9078"
    "D10" = "0.1412"
    "D11" = "0.07370"
    "D12" = "0.03355"
    "D14" = "0.09334"
    "D15" = "3.847"
    "D16" = "0.001590"
    "D17" = "0.04647"
    "D18" = "0.0005938"
    "D19" = "0.006084"
    "E20" = "19HotbitTokenHTBBestin24h"
    "D21" = "0.0009847"
    "D22" = "0.00007798"
    "D23" = "0.0002899"
    "D24" = "3.614"
    "D27" = "0.1293"
    "D40" = "0.03885"
    "D41" = "0.006230"
    "E41" = "40KickTokenKICK"
    "D42" = "0.1071"
    "D43" = "0.002799"
    "D44" = "0.007208"
    "D45" = "0.00005184"
    "D50" = "0.00002099"
    "D51" = "0.0001999"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
